$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying records (sightings) for rows 4-12 were re-matched to
# coordinates: each row keeps its location/date/reporter columns, but the
# species-identification columns (A,B,D,E,F,G,H), the coordinate columns
# (Q,R, now rounded to whole metres) and a couple of occurrence-specific
# columns (K/L/M/N blank markers, AC public comment) move between rows.
# The Starttid/Sluttid columns (Z, AB) are dropped for every data row.

$rows = @{
    4  = @{ A = 111936768; B = 90087; D = "LC"; E = 3298; F = "Trådticka"; G = "Climacocystis borealis"; H = "(Fr.) Kotl. & Pouzar"; Q = 490317; R = 7088522; KLMN = $false; AC = $null }
    5  = @{ A = 111936777; B = 77515; D = "NT"; E = 6425; F = "Garnlav"; G = "Alectoria sarmentosa"; H = "(Ach.) Ach."; Q = 490056; R = 7088709; KLMN = $false; AC = $null }
    6  = @{ A = 111936775; B = 89419; D = "NT"; E = 1204; F = "Gränsticka"; G = "Phellopilus nigrolimitatus"; H = "(Romell) Niemelä, T.Wagner & M.Fisch."; Q = 490380; R = 7088379; KLMN = $false; AC = $null }
    7  = @{ A = 111936776; B = 77515; D = "NT"; E = 6425; F = "Garnlav"; G = "Alectoria sarmentosa"; H = "(Ach.) Ach."; Q = 490398; R = 7088445; KLMN = $false; AC = $null }
    8  = @{ A = 111936779; B = 77515; D = "NT"; E = 6425; F = "Garnlav"; G = "Alectoria sarmentosa"; H = "(Ach.) Ach."; Q = 490008; R = 7088597; KLMN = $false; AC = $null }
    9  = @{ A = 111936780; B = 77515; D = "NT"; E = 6425; F = "Garnlav"; G = "Alectoria sarmentosa"; H = "(Ach.) Ach."; Q = 489952; R = 7088557; KLMN = $false; AC = $null }
    10 = @{ A = 111936767; B = 90087; D = "LC"; E = 3298; F = "Trådticka"; G = "Climacocystis borealis"; H = "(Fr.) Kotl. & Pouzar"; Q = 490377; R = 7088412; KLMN = $false; AC = $null }
    11 = @{ A = 111936781; B = 89793; D = "LC"; E = 4217; F = "Blodticka"; G = "Meruliopsis taxicola"; H = "(Pers.:Fr.) Bondartsev"; Q = 490315; R = 7088552; KLMN = $false; AC = $null }
    12 = @{ A = 111936774; B = 56414; D = "NT"; E = 100049; F = "Spillkråka"; G = "Dryocopus martius"; H = "(Linnaeus, 1758)"; Q = 490378; R = 7088551; KLMN = $true; AC = "hack" }
}

foreach ($r in 4..12) {
    $rec = $rows[$r]

    $ws.Range("A$r").Value = $rec.A
    $ws.Range("B$r").Value = $rec.B
    $ws.Range("D$r").Value = $rec.D
    $ws.Range("E$r").Value = $rec.E
    $ws.Range("F$r").Value = $rec.F
    $ws.Range("G$r").Value = $rec.G
    $ws.Range("H$r").Value = $rec.H
    $ws.Range("Q$r").Value = $rec.Q
    $ws.Range("R$r").Value = $rec.R

    # Starttid/Sluttid no longer present on any of these rows.
    $ws.Range("Z$r").ClearContents()
    $ws.Range("AB$r").ClearContents()

    # Ålder-Stadium / Kön / Aktivitet / Metod — blank marker cells that
    # only exist on whichever row currently carries the Spillkråka record.
    if ($rec.KLMN) {
        $ws.Range("K$r").Value = ""
        $ws.Range("L$r").Value = ""
        $ws.Range("M$r").Value = ""
        $ws.Range("N$r").Value = ""
    } else {
        $ws.Range("K$r").ClearContents()
        $ws.Range("L$r").ClearContents()
        $ws.Range("M$r").ClearContents()
        $ws.Range("N$r").ClearContents()
    }

    # Publik kommentar — only set on the row that currently carries the
    # "hack" observation.
    if ($rec.AC) {
        $ws.Range("AC$r").Value = $rec.AC
    } else {
        $ws.Range("AC$r").ClearContents()
    }
}
